# Edit script for 09-A-ProjectWork.docx - "Fixed typos in 09"
$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 4 (process first, bottom of document, to keep earlier offsets stable):
# Remove manual page-break paragraph and the lastRenderedPageBreak hint
# before "When one team member has pushed a commit..."
# -----------------------------------------------------------------
$pBreakPara = $d.Paragraphs(89)
$rBreak = $pBreakPara.Range
$delBreak = $d.Range($rBreak.Start, $rBreak.Start + 1)
$delBreak.Text = ""

$rngWhen = $d.Content
$needleWhen = "When one team member has pushed a commit to the feature branch on the origin the team should make a Draft Pull Request"
$foundWhen = $rngWhen.Find.Execute($needleWhen, $false, $false, $false, $false, $false, $true, 1, $false, $needleWhen, 2)

# -----------------------------------------------------------------
# Edit 3: "git pull origin branchname" -> "git fetch origin <branchname>"
# plus a new paragraph "git switch <branchname>" right after it.
# -----------------------------------------------------------------
$pGitPull = $d.Paragraphs(82)
$rGitPull = $pGitPull.Range
# Insert the new empty paragraph FIRST (before editing text) so it
# inherits paragraph 82's original (Courier / indented) formatting.
$rGitPull.InsertParagraphAfter()

$pSwitch = $d.Paragraphs(83)
$rSwitch = $pSwitch.Range
$insSwitch = $d.Range($rSwitch.Start, $rSwitch.Start)
$insSwitch.Text = "git switch <branchname>"

$pGitPull2 = $d.Paragraphs(82)
$rGitPull2 = $pGitPull2.Range
$fullGitPull = $d.Range($rGitPull2.Start, $rGitPull2.End - 1)
$needleGitPull = "git pull origin branchname"
$foundGitPull = $fullGitPull.Find.Execute($needleGitPull, $false, $false, $false, $false, $false, $true, 1, $false, "git fetch origin <branchname>", 2)

# -----------------------------------------------------------------
# Edit 2: "Have all team members pull the new feature branch into
# their local repositories." ->
# "Have all team members *fetch* a copy of the new feature branch
# into their local repositories." ("fetch" italic)
# -----------------------------------------------------------------
$pHaveAll = $d.Paragraphs(81)
$rHaveAll = $pHaveAll.Range
$baseHaveAll = $rHaveAll.Start
$pullStart = $baseHaveAll + 22
$pullEnd = $pullStart + 4

$insCopyOf = $d.Range($pullEnd, $pullEnd)
$insCopyOf.Text = " a copy of"

$pullRange = $d.Range($pullStart, $pullEnd)
$pullRange.Text = "fetch"
$pullRange.Font.Italic = $true

# -----------------------------------------------------------------
# Edit 1: "git set remote upstream" -> "git remote add upstream"
# -----------------------------------------------------------------
$rngUpstream = $d.Content
$foundUpstream = $rngUpstream.Find.Execute("git set remote upstream", $false, $false, $false, $false, $false, $true, 1, $false, "git remote add upstream", 2)

Write-Output "foundWhen=$foundWhen foundGitPull=$foundGitPull foundUpstream=$foundUpstream"
